# Weekly price update: a new week's record is inserted at the top of the
# data block (row 94), pushing every existing record down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 94 - this shifts rows 94:186 down to
# 95:187 and extends the used range to A1:R187 automatically.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(94, 1).Value = 7
$ws.Cells.Item(94, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(94, 3).Value = "Ñuble"
$ws.Cells.Item(94, 4).Value = "2022-01-28"
$ws.Cells.Item(94, 5).Value = 16
$ws.Cells.Item(94, 6).Value = 100112032
$ws.Cells.Item(94, 7).Value = "Zapallo italiano"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 120
$ws.Cells.Item(94, 11).Value = 8500
$ws.Cells.Item(94, 12).Value = 9000
$ws.Cells.Item(94, 13).Value = 8750
$ws.Cells.Item(94, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(94, 15).Value = "Región del Maule"
$ws.Cells.Item(94, 16).Value = 146
$ws.Cells.Item(94, 17).Value = 60
$ws.Cells.Item(94, 18).Value = "Hortaliza"
